$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38..151 down to 39..152
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with the new data record
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44414
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112001
$ws.Range("G38").Value = "Berenjena"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 14000
$ws.Range("L38").Value = 14000
$ws.Range("M38").Value = 14000
$ws.Range("N38").Value = "$/caja 60 unidades"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 233
$ws.Range("Q38").Value = 60
$ws.Range("R38").Value = "Hortaliza"
